$d = $word.ActiveDocument

# Locate the paragraph that contains the "2200080 : Avoid using Thread.Sleep ..."
# release-note entry so we don't depend on a fixed paragraph index.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "2200080") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Delete the "2200080 : Avoid using Thread.Sleep with dynamic parameter in a
    # controller action" paragraph entirely (including its paragraph mark).
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.Delete()

    # The paragraph immediately preceding it is the blank spacer paragraph that
    # separated it from the previous release-note entry; remove that too so the
    # remaining entries keep a single blank line between them.
    $prevIndex = $targetIndex - 1
    if ($prevIndex -ge 1) {
        $prev = $d.Paragraphs.Item($prevIndex)
        if ($prev.Range.Text.Trim() -eq "") {
            $prev.Range.Delete()
        }
    }
}
